$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3655.4546
$ws.Range("J17").Value = 3655.4546
$ws.Range("L17").Value = 10966.3638
$ws.Range("N17").Value = -11302.3638
$ws.Range("H40").Value = 2054.8708
$ws.Range("I40").Value = 1859.1364
$ws.Range("K40").Value = 1859.1364
$ws.Range("M40").Value = -1684.1364
$ws.Range("H70").Value = 1408
$ws.Range("I70").Value = 1634
$ws.Range("J70").Value = 1238.5
$ws.Range("K70").Value = 4902
$ws.Range("L70").Value = 3715.5
$ws.Range("M70").Value = -4632
$ws.Range("N70").Value = -4255.5
$ws.Range("H73").Value = 1408
$ws.Range("I73").Value = 1634
$ws.Range("J73").Value = 1238.5
$ws.Range("K73").Value = 4902
$ws.Range("L73").Value = 3715.5
$ws.Range("M73").Value = -3966
$ws.Range("N73").Value = -5587.5
$ws.Range("H137").Value = 1141.129
$ws.Range("I137").Value = 919.1395
$ws.Range("J137").Value = 1643.5264
$ws.Range("K137").Value = 2757.4185
$ws.Range("L137").Value = 4930.5792
$ws.Range("M137").Value = -207.4184999999998
$ws.Range("N137").Value = -10030.5792
$ws.Range("H138").Value = 469689.12
$ws.Range("I138").Value = 1661.2106
$ws.Range("J138").Value = 589858.4399999999
$ws.Range("K138").Value = 4983.6318
$ws.Range("L138").Value = 1769575.32
$ws.Range("M138").Value = 156.3681999999999
$ws.Range("N138").Value = -1779855.32

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7293.297
$ws.Range("I32").Value = 6022.533
$ws.Range("K32").Value = 6022.533
$ws.Range("M32").Value = -5735.533
$ws.Range("H61").Value = 111112536
$ws.Range("I61").Value = 142858400
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 142858400
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -142858188
$ws.Range("N61").Value = -2424
$ws.Range("H74").Value = 2016.25
$ws.Range("I74").Value = 1294.8823
$ws.Range("K74").Value = 1294.8823
$ws.Range("M74").Value = -420.8823
$ws.Range("H77").Value = 2016.25
$ws.Range("I77").Value = 1294.8823
$ws.Range("K77").Value = 6474.4115
$ws.Range("M77").Value = -2106.4115
$ws.Range("H102").Value = 8334487.5
$ws.Range("I102").Value = 11905790
$ws.Range("J102").Value = 1448
$ws.Range("K102").Value = 11905790
$ws.Range("L102").Value = 1448
$ws.Range("M102").Value = -11904168
$ws.Range("N102").Value = -4692
$ws.Range("H132").Value = 2517.3096
$ws.Range("I132").Value = 2015.32
$ws.Range("J132").Value = 3255.5293
$ws.Range("K132").Value = 6045.96
$ws.Range("L132").Value = 9766.5879
$ws.Range("M132").Value = -3515.96
$ws.Range("N132").Value = -14826.5879
$ws.Range("H136").Value = 111112536
$ws.Range("I136").Value = 142858400
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 428575200
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -428572650
$ws.Range("N136").Value = -11100

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2550.4856
$ws.Range("I86").Value = 2489.5
$ws.Range("J86").Value = 2653.6924
$ws.Range("K86").Value = 2489.5
$ws.Range("L86").Value = 2653.6924
$ws.Range("M86").Value = -1366.5
$ws.Range("N86").Value = -4899.6924
$ws.Range("H89").Value = 2550.4856
$ws.Range("I89").Value = 2489.5
$ws.Range("J89").Value = 2653.6924
$ws.Range("K89").Value = 12447.5
$ws.Range("L89").Value = 13268.462
$ws.Range("M89").Value = -6831.5
$ws.Range("N89").Value = -24500.462
$ws.Range("H134").Value = 1314
$ws.Range("I134").Value = 1116.3334
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 3349.0002
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -814.0001999999999
$ws.Range("N134").Value = -12570

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 71430030
$ws.Range("I16").Value = 111112440
$ws.Range("J16").Value = 1710
$ws.Range("K16").Value = 111112440
$ws.Range("L16").Value = 1710
$ws.Range("M16").Value = -111112153
$ws.Range("N16").Value = -2284
$ws.Range("H31").Value = 1357.5853
$ws.Range("I31").Value = 1357.5853
$ws.Range("K31").Value = 1357.5853
$ws.Range("M31").Value = -1062.5853
$ws.Range("H34").Value = 1357.5853
$ws.Range("I34").Value = 1357.5853
$ws.Range("K34").Value = 1357.5853
$ws.Range("M34").Value = -1155.5853
$ws.Range("H58").Value = 4381.5
$ws.Range("I58").Value = 1243.2632
$ws.Range("J58").Value = 7519.737
$ws.Range("K58").Value = 1243.2632
$ws.Range("L58").Value = 7519.737
$ws.Range("M58").Value = -1040.2632
$ws.Range("N58").Value = -7925.737
$ws.Range("H105").Value = 810
$ws.Range("I105").Value = 779.8
$ws.Range("K105").Value = 779.8
$ws.Range("M105").Value = 967.2
$ws.Range("H113").Value = 71430030
$ws.Range("I113").Value = 111112440
$ws.Range("J113").Value = 1710
$ws.Range("K113").Value = 111112440
$ws.Range("L113").Value = 1710
$ws.Range("M113").Value = -111110270
$ws.Range("N113").Value = -6050
$ws.Range("H136").Value = 4381.5
$ws.Range("I136").Value = 1243.2632
$ws.Range("J136").Value = 7519.737
$ws.Range("K136").Value = 3729.7896
$ws.Range("L136").Value = 22559.211
$ws.Range("M136").Value = -1179.7896
$ws.Range("N136").Value = -27659.211

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3907.75
$ws.Range("J39").Value = 3923.1428
$ws.Range("L39").Value = 11769.4284
$ws.Range("N39").Value = -12357.4284
$ws.Range("H46").Value = 2099.6667
$ws.Range("I46").Value = 2649.5
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 7948.5
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = -7857.5
$ws.Range("N46").Value = -3182
$ws.Range("H55").Value = 1725.2354
$ws.Range("J55").Value = 1725.2354
$ws.Range("L55").Value = 5175.706200000001
$ws.Range("N55").Value = -5529.706200000001
$ws.Range("H113").Value = 684.6667
$ws.Range("I113").Value = 586
$ws.Range("J113").Value = 727.56525
$ws.Range("K113").Value = 1758
$ws.Range("L113").Value = 2182.69575
$ws.Range("M113").Value = 412
$ws.Range("N113").Value = -6522.69575
$ws.Range("H131").Value = 20836216
$ws.Range("J131").Value = 3596.1892
$ws.Range("L131").Value = 10788.5676
$ws.Range("N131").Value = -20868.5676

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 28127846
$ws.Range("I70").Value = 27780596
$ws.Range("J70").Value = 28574312
$ws.Range("K70").Value = 27780596
$ws.Range("L70").Value = 28574312
$ws.Range("M70").Value = -27780326
$ws.Range("N70").Value = -28574852
$ws.Range("H73").Value = 28127846
$ws.Range("I73").Value = 27780596
$ws.Range("J73").Value = 28574312
$ws.Range("K73").Value = 27780596
$ws.Range("L73").Value = 28574312
$ws.Range("M73").Value = -27779660
$ws.Range("N73").Value = -28576184
$ws.Range("H80").Value = 5281.8184
$ws.Range("I80").Value = 5440
$ws.Range("K80").Value = 5440
$ws.Range("M80").Value = -4442
$ws.Range("H83").Value = 5281.8184
$ws.Range("I83").Value = 5440
$ws.Range("K83").Value = 27200
$ws.Range("M83").Value = -22208
$ws.Range("H132").Value = 2675.2285
$ws.Range("I132").Value = 2564.5908
$ws.Range("K132").Value = 7693.7724
$ws.Range("M132").Value = -5163.7724

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3008.2222
$ws.Range("I7").Value = 2710.6667
$ws.Range("J7").Value = 3603.3333
$ws.Range("K7").Value = 2710.6667
$ws.Range("L7").Value = 3603.3333
$ws.Range("M7").Value = -2598.6667
$ws.Range("N7").Value = -3827.3333
$ws.Range("H123").Value = 40960
$ws.Range("J123").Value = 40960
$ws.Range("L123").Value = 40960
$ws.Range("N123").Value = -50760
$ws.Range("H126").Value = 3008.2222
$ws.Range("I126").Value = 2710.6667
$ws.Range("J126").Value = 3603.3333
$ws.Range("K126").Value = 8132.000100000001
$ws.Range("L126").Value = 10809.9999
$ws.Range("M126").Value = -5662.000100000001
$ws.Range("N126").Value = -15749.9999
$ws.Range("H136").Value = 2700
$ws.Range("I136").Value = 2000
$ws.Range("J136").Value = 2933.3333
$ws.Range("K136").Value = 6000
$ws.Range("L136").Value = 8799.999899999999
$ws.Range("M136").Value = -3450
$ws.Range("N136").Value = -13899.9999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 5020
$ws.Range("I43").Value = 5020
$ws.Range("K43").Value = 5020
$ws.Range("M43").Value = -4871
$ws.Range("H132").Value = 3187.389
$ws.Range("I132").Value = 3433.6365
$ws.Range("J132").Value = 2800.4285
$ws.Range("K132").Value = 10300.9095
$ws.Range("L132").Value = 8401.2855
$ws.Range("M132").Value = -7770.9095
$ws.Range("N132").Value = -13461.2855
$ws.Range("H136").Value = 1924.6154
$ws.Range("I136").Value = 1852.5
$ws.Range("J136").Value = 1986.4286
$ws.Range("K136").Value = 5557.5
$ws.Range("L136").Value = 5959.2858
$ws.Range("M136").Value = -3007.5
$ws.Range("N136").Value = -11059.2858
